$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5234.4
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 5543
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 5543
$ws.Range("M32").Value = -3674
$ws.Range("N32").Value = -6195

$ws.Range("H82").Value = 220.5
$ws.Range("I82").Value = 220.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 661.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -255.5
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 220.5
$ws.Range("I85").Value = 220.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 661.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 742.5
$ws.Range("N85").ClearContents()

$ws.Range("H137").Value = 4768.4287
$ws.Range("I137").Value = 3855.2
$ws.Range("K137").Value = 11565.6
$ws.Range("M137").Value = -9015.599999999999

$ws.Range("H138").Value = 6541.6113
$ws.Range("I138").Value = 4124.25
$ws.Range("J138").Value = 7232.2856
$ws.Range("K138").Value = 12372.75
$ws.Range("L138").Value = 21696.8568
$ws.Range("M138").Value = -7232.75
$ws.Range("N138").Value = -31976.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1499.75
$ws.Range("I61").Value = 1499.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1499.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1287.75
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 1993.08
$ws.Range("I74").Value = 1811.6842
$ws.Range("J74").Value = 2567.5
$ws.Range("K74").Value = 1811.6842
$ws.Range("L74").Value = 2567.5
$ws.Range("M74").Value = -937.6841999999999
$ws.Range("N74").Value = -4315.5

$ws.Range("H77").Value = 1993.08
$ws.Range("I77").Value = 1811.6842
$ws.Range("J77").Value = 2567.5
$ws.Range("K77").Value = 9058.421
$ws.Range("L77").Value = 12837.5
$ws.Range("M77").Value = -4690.421
$ws.Range("N77").Value = -21573.5

$ws.Range("H113").Value = 52398
$ws.Range("J113").Value = 52398
$ws.Range("L113").Value = 52398
$ws.Range("N113").Value = -61076

$ws.Range("H114").Value = 19996.666
$ws.Range("J114").Value = 19996.666
$ws.Range("L114").Value = 19996.666
$ws.Range("N114").Value = -28674.666

$ws.Range("H132").Value = 3024.1765
$ws.Range("I132").Value = 3286.5715
$ws.Range("J132").Value = 1799.6666
$ws.Range("K132").Value = 9859.7145
$ws.Range("L132").Value = 5398.9998
$ws.Range("M132").Value = -7329.7145
$ws.Range("N132").Value = -10458.9998

$ws.Range("H136").Value = 1499.75
$ws.Range("I136").Value = 1499.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4499.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1949.25
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1103.5
$ws.Range("I14").Value = 1103.5
$ws.Range("K14").Value = 3310.5
$ws.Range("M14").Value = -3137.5

$ws.Range("H55").Value = 3599.8
$ws.Range("I55").Value = 2500
$ws.Range("J55").Value = 3874.75
$ws.Range("K55").Value = 7500
$ws.Range("L55").Value = 11624.25
$ws.Range("M55").Value = -7323
$ws.Range("N55").Value = -11978.25

$ws.Range("H68").Value = 5940
$ws.Range("J68").Value = 8192
$ws.Range("L68").Value = 24576
$ws.Range("N68").Value = -26198

$ws.Range("H71").Value = 5940
$ws.Range("J71").Value = 8192
$ws.Range("L71").Value = 73728
$ws.Range("N71").Value = -81840

$ws.Range("H86").Value = 1496.3334
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1496.3334
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4489.0002
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -6861.0002

$ws.Range("H89").Value = 1496.3334
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1496.3334
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 13467.0006
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -25323.0006

$ws.Range("H122").Value = 2316.6667
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2475
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 22275
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = -27175

$ws.Range("H124").Value = 4996.75
$ws.Range("I124").Value = 4996.75
$ws.Range("K124").Value = 14990.25
$ws.Range("M124").Value = -10080.25

$ws.Range("H130").Value = 4232.5
$ws.Range("I130").Value = 4232.5
$ws.Range("K130").Value = 12697.5
$ws.Range("M130").Value = -7677.5

$ws.Range("H131").Value = 1979.909
$ws.Range("I131").Value = 1908.7778
$ws.Range("J131").Value = 2300
$ws.Range("K131").Value = 5726.3334
$ws.Range("L131").Value = 6900
$ws.Range("M131").Value = -686.3334000000004
$ws.Range("N131").Value = -16980

$ws.Range("H133").Value = 9822.200000000001
$ws.Range("I133").Value = 9822.200000000001
$ws.Range("K133").Value = 29466.6
$ws.Range("M133").Value = -24406.6

$ws.Range("H134").Value = 7521.8667
$ws.Range("I134").Value = 913.9286
$ws.Range("K134").Value = 2741.7858
$ws.Range("M134").Value = 2328.2142

$ws.Range("H140").Value = 592010.8
$ws.Range("I140").Value = 836186.5600000001
$ws.Range("J140").Value = 5989
$ws.Range("K140").Value = 2508559.68
$ws.Range("L140").Value = 17967
$ws.Range("M140").Value = -2503379.68
$ws.Range("N140").Value = -28327

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6174.8
$ws.Range("I7").Value = 7542.3335
$ws.Range("J7").Value = 4123.5
$ws.Range("K7").Value = 7542.3335
$ws.Range("L7").Value = 4123.5
$ws.Range("M7").Value = -7430.3335
$ws.Range("N7").Value = -4347.5

$ws.Range("H40").Value = 5575.8237
$ws.Range("I40").Value = 3056.1428
$ws.Range("J40").Value = 17334.334
$ws.Range("K40").Value = 3056.1428
$ws.Range("L40").Value = 17334.334
$ws.Range("M40").Value = -2920.1428
$ws.Range("N40").Value = -17606.334

$ws.Range("H122").Value = 10638.044
$ws.Range("I122").Value = 8272.684999999999
$ws.Range("J122").Value = 21873.5
$ws.Range("K122").Value = 24818.055
$ws.Range("L122").Value = 65620.5
$ws.Range("M122").Value = -22368.055
$ws.Range("N122").Value = -70520.5

$ws.Range("H126").Value = 6174.8
$ws.Range("I126").Value = 7542.3335
$ws.Range("J126").Value = 4123.5
$ws.Range("K126").Value = 22627.0005
$ws.Range("L126").Value = 12370.5
$ws.Range("M126").Value = -20157.0005
$ws.Range("N126").Value = -17310.5
